$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:B4")

# Force text format first so Excel doesn't auto-convert the date-looking
# and numeric-looking strings into real dates/numbers (which would change
# the underlying cell value away from the literal text in the diff).
$rng.NumberFormat = "@"

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Total Amount"

$ws.Range("A2").Value = "2023-06-18"
$ws.Range("B2").Value = "2387.00"

$ws.Range("A3").Value = "2023-06-06"
$ws.Range("B3").Value = "213.00"

$ws.Range("A4").Value = "2023-06-15"
$ws.Range("B4").Value = "50.00"

# Restore the cells' formatting to the default (no explicit style index),
# matching the original workbook where these cells had no style applied.
$rng.ClearFormats()
